$wb = $excel.ActiveWorkbook

# --- Sheet "df_devolucoes": add row 4 with perc_sim = 49.83 ---
$wsDev = $wb.Worksheets.Item("df_devolucoes")
$wsDev.Range("A4").Value = "perc_sim"
$wsDev.Range("B4").Value = 49.83

# --- Sheet "df_produtos": reorder/relabel product names in column A ---
$wsProd = $wb.Worksheets.Item("df_produtos")
$wsProd.Range("A2").Value = "SPLIT 10001 À 18000"
$wsProd.Range("A3").Value = "SPLIT 18001 À 30000"
$wsProd.Range("A4").Value = "JANELA ATÉ 8500 BTUS"
$wsProd.Range("A5").Value = "TV LED DE 48`" À 54`""
$wsProd.Range("A6").Value = "GRILL"
$wsProd.Range("A7").Value = "SANDUICHEIRA"
$wsProd.Range("A8").Value = "SMART MULTISIM LIVRE"
$wsProd.Range("A9").Value = "CAFETEIRA"
$wsProd.Range("A10").Value = "TV 4K DE 60`" À 69`""
$wsProd.Range("A11").Value = "TV 4K ATÉ 59`""
